$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Naive component forecaster bug fix: C2 was an erroneous stray value - remove it entirely.
$ws.Range("C2").ClearContents()

# Recalculated forecast values (tiny floating point corrections from the fix).
$ws.Range("E2").Value = 6.253707197847591
$ws.Range("C3").Value = -7.921319741078636
$ws.Range("E3").Value = -14.35806537048446
$ws.Range("E5").Value = 6.136355062499965
$ws.Range("C6").Value = 1.477633171193116
$ws.Range("C7").Value = -0.2005250704869121
$ws.Range("C8").Value = 2.234527904461148
$ws.Range("C9").Value = 1.311727872618218
$ws.Range("E9").Value = 1.093673275363694
$ws.Range("C10").Value = 1.784808447869191
$ws.Range("C12").Value = 2.159589514946725
$ws.Range("C13").Value = 0.8014493436638848
$ws.Range("E13").Value = 1.609625625599986
$ws.Range("C14").Value = -3.107661574595766
$ws.Range("E14").Value = -8.513835774400015
$ws.Range("C15").Value = -1.621578487659103
$ws.Range("C16").Value = 1.906376895025041
$ws.Range("E16").Value = 0.9449384537270955
$ws.Range("C17").Value = -0.4925007786849234
